$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "달력 만들기"
$ws.Range("C13").Value = "달력 관리"
$ws.Range("C14").Value = "멤버 관리"

$ws.Range("C13").Select()
